$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns to match the refreshed crypto feed.
# D-column values are forced to text (NumberFormat "@") before assignment so that
# numeric-looking prices (e.g. "6.60", "0.997") are preserved exactly as strings
# instead of being coerced into doubles (which would lose trailing zeros / introduce
# floating point noise). Style is reset to Normal afterwards so no stray formatting
# is left behind on the cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.057.53"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.86%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.651.68"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.44%  "

$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "533.21"
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.21"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.18%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.13%  "

$ws.Range("E8").Value = "  +0.75%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.60"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.58%  "

$ws.Range("E10").Value = "  +5.34%  "

$ws.Range("E11").Value = "  +1.75%  "

$ws.Range("E12").Value = "  +0.02%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.114.10"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.90%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "61.033.56"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.83%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "22.08"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.24%  "

$ws.Range("E16").Value = "  +2.56%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.660.96"
$ws.Range("D17").Style = "Normal"

$ws.Range("E18").Value = "  +0.20%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "355.30"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.34%  "

$ws.Range("E20").Value = "  +0.67%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.26"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.61%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.11%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "61.65"
$ws.Range("D23").Style = "Normal"

$ws.Range("E24").Value = "  +2.03%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.168"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.39%  "

$ws.Range("E26").Value = "  +0.21%  "

$ws.Range("E27").Value = "  +2.11%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.40"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.44%  "

$ws.Range("E29").Value = "  -0.14%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.21"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.06%  "

$ws.Range("E31").Value = "  +4.33%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.57"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.69%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "149.99"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.29%  "

$ws.Range("E34").Value = "  +3.57%  "

$ws.Range("E35").Value = "  +1.23%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.921"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +8.77%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.894"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.13%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "310.15"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.79%  "

$ws.Range("E39").Value = "  +0.83%  "

$ws.Range("E40").Value = "  +1.50%  "

$ws.Range("E41").Value = "  +3.43%  "

$ws.Range("E42").Value = "  +1.74%  "

$ws.Range("E43").Value = "  +1.64%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.997"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.04%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.98"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.66%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.97"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.27%  "

$ws.Range("E47").Value = "  +2.43%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "19.24"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +8.16%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.35"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.31%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.995.52"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.34%  "

$ws.Range("E51").Value = "  +2.72%  "
